$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.040.60"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.649.70"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.71"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5295"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2606"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06297"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.28"
$ws.Range("E10").Value = "  -3.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07745"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.472"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "1.653.58"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5439"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "0.0₅8091"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.05"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "26.070.57"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.548"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.39"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.02"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.986"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "139.72"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1242"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.245"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.437"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05915"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.502"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.237"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.547"
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.414"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9402"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.755"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5678"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01604"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.846"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8438"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "1.008.88"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.57"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "1.799.41"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.86"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4298"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.476"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.857"
$ws.Range("E50").Value = "  -3.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05153"
$ws.Range("E51").Value = "  -0.52%  "
